# Weekly data refresh: a new week's price observation is inserted as a new
# row 65 on Sheet1 ("Hortaliza, Vega Central Mapocho de Santiago - Jengibre").
# Inserting the row pushes the previously existing rows 65-130 down to
# rows 66-131 (and grows the sheet's used range from A1:R130 to A1:R131),
# matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65 - everything below shifts down by one.
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A65").Value = 9
$ws.Range("B65").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C65").Value = "Metropolitana"
$ws.Range("D65").Value = 45068
$ws.Range("E65").Value = 13
$ws.Range("F65").Value = 100114007
$ws.Range("G65").Value = "Jengibre"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 520
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 16000
$ws.Range("M65").Value = 15481
$ws.Range("N65").Value = "`$/caja 13 kilos"
$ws.Range("O65").Value = "Perú"
$ws.Range("P65").Value = 1191
$ws.Range("Q65").Value = 13
$ws.Range("R65").Value = "Hortaliza"
